$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 23:40"

# Row 4
$ws.Range("B4").Value = 7592465
$ws.Range("C4").Value = 42430
$ws.Range("D4").Value = 4808668
$ws.Range("E4").Value = 2569705
$ws.Range("G4").Value = 570
$ws.Range("H4").Value = 214092

# Row 6
$ws.Range("B6").Value = 4906760
$ws.Range("C6").Value = 24529
$ws.Range("E6").Value = 528182
$ws.Range("G6").Value = 554
$ws.Range("H6").Value = 145985

# Row 15
$ws.Range("A15").Value = "Reino Unido"
$ws.Range("B15").Value = 480017
$ws.Range("C15").Value = 7070
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 49
$ws.Range("H15").Value = 42317

# Row 16
$ws.Range("A16").Value = "Chile"
$ws.Range("B16").Value = 468471
$ws.Range("C16").Value = 1881
$ws.Range("D16").Value = 440881
$ws.Range("E16").Value = 14671
$ws.Range("G16").Value = 52
$ws.Range("H16").Value = 12919

# Row 17
$ws.Range("A17").Value = "Iran"
$ws.Range("B17").Value = 468119
$ws.Range("C17").Value = 3523
$ws.Range("D17").Value = 387675
$ws.Range("E17").Value = 53698
$ws.Range("G17").Value = 179
$ws.Range("H17").Value = 26746

# Row 25
$ws.Range("B25").Value = 300028
$ws.Range("C25").Value = 1665
$ws.Range("E25").Value = 30931

# Row 30
$ws.Range("B30").Value = 140351
$ws.Range("C30").Value = 817
$ws.Range("E30").Value = 16458
$ws.Range("G30").Value = 102
$ws.Range("H30").Value = 11597

# Row 52
$ws.Range("A52").Value = "Costa Rica"
$ws.Range("B52").Value = 79182
$ws.Range("C52").Value = 1353
$ws.Range("D52").Value = 45007
$ws.Range("E52").Value = 33225
$ws.Range("G52").Value = 20
$ws.Range("H52").Value = 950

# Row 53
$ws.Range("A53").Value = "Honduras"
$ws.Range("B53").Value = 78269
$ws.Range("C53").Value = 671
$ws.Range("D53").Value = 28978
$ws.Range("E53").Value = 46905
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 2386

# Row 54
$ws.Range("A54").Value = "Portugal"
$ws.Range("B54").Value = 78247
$ws.Range("C54").Value = 963
$ws.Range("D54").Value = 49845
$ws.Range("E54").Value = 26407
$ws.Range("G54").Value = 12
$ws.Range("H54").Value = 1995

# Row 55
$ws.Range("A55").Value = "Etiopia"
$ws.Range("B55").Value = 77860
$ws.Range("C55").Value = 872
$ws.Range("D55").Value = 32325
$ws.Range("E55").Value = 44321
$ws.Range("G55").Value = 6
$ws.Range("H55").Value = 1214

# Row 59
$ws.Range("B59").Value = 58238
$ws.Range("C59").Value = 784
$ws.Range("D59").Value = 54854
$ws.Range("E59").Value = 2907
$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 477

# Row 63
$ws.Range("A63").Value = "Argelia"
$ws.Range("B63").Value = 51995
$ws.Range("C63").Value = 148
$ws.Range("D63").Value = 36482
$ws.Range("E63").Value = 13757
$ws.Range("G63").Value = 7
$ws.Range("H63").Value = 1756

# Row 64
$ws.Range("A64").Value = "Armenia"
$ws.Range("B64").Value = 51925
$ws.Range("C64").Value = 543
$ws.Range("D64").Value = 44583
$ws.Range("E64").Value = 6370
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 972

# Row 124
$ws.Range("B124").Value = 5370
$ws.Range("C124").Value = 159
$ws.Range("D124").Value = 2436
$ws.Range("E124").Value = 2741
$ws.Range("G124").Value = 4
$ws.Range("H124").Value = 193

# Row 189
$ws.Range("B189").Value = 222
$ws.Range("C189").Value = 1
$ws.Range("E189").Value = 32

# Row 215
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# Row 216
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
